# Add a Calibri w:rFonts entry to every table run that already carries
# explicit run formatting (w:b / w:i / w:strike / w:color / w:sz).
#
# This Find/Replace in this runtime operates over the whole document
# (not scoped to the Range it was obtained from), so every search string
# used below is unique in the document -- except "Total" which
# legitimately appears twice in the table and needs the exact same
# rFonts applied both times, so a single global replace handles both.

$d = $word.ActiveDocument

function Set-CalibriByText($searchText) {
    $rng = $d.Content
    $find = $rng.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Replacement.Font.NameAscii = "Calibri"
    $find.Replacement.Font.NameFarEast = "Calibri"
    $find.Replacement.Font.NameOther = "Calibri"
    $find.Replacement.Font.NameBi = "Calibri"
    $find.Text = $searchText
    $find.Replacement.Text = $searchText
    $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)
}

# Cells with real text content (14 cells; "Total" hits 2 cells at once).
$texts = @(
    "Item", "Price", "Quantity", "Total",
    "Apple", "1,76 €", "23", "40,48 €",
    "Banana", "1,99 €", "45", "89,55 €",
    "130,03 €"
)
foreach ($txt in $texts) {
    Set-CalibriByText $txt
}

# The two formatted-but-empty cells (row 3, cols 2 and 4) have no text to
# search for, so temporarily give them unique marker text, format that
# marker, then strip the marker text back out while keeping the (now
# formatted) empty run in place.
$t = $d.Tables.Item(1)
$cell32 = $t.Cell(3, 2)
$cell32.Range.InsertBefore("ZZTEMPMARKERA")
$t = $d.Tables.Item(1)
$cell34 = $t.Cell(3, 4)
$cell34.Range.InsertBefore("ZZTEMPMARKERB")

Set-CalibriByText "ZZTEMPMARKERA"
Set-CalibriByText "ZZTEMPMARKERB"

$t = $d.Tables.Item(1)
$cell32 = $t.Cell(3, 2)
$r32 = $cell32.Range
$r32.MoveEnd(1, -2)
$r32.Text = ""

$t = $d.Tables.Item(1)
$cell34 = $t.Cell(3, 4)
$r34 = $cell34.Range
$r34.MoveEnd(1, -2)
$r34.Text = ""

Write-Host "done"
